$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts C:Z -> D:AA)
$ws.Columns("C:C").Insert()

# New "pid" header + values
$ws.Range("C1").Value = "pid"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5

# Match the author's final selection/view state
$ws.Range("C7").Select()
